$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wrapper_ready")
$ws.Activate()

$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 2
$ws.Range("L5").Value = 2

$ws.Range("L6").Select()
